$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add the new "FiP" (5th Pref) worksheet right after the existing "FoP"
#    sheet, i.e. at the very end of the workbook.
# ---------------------------------------------------------------------------
$fop = $wb.Worksheets.Item("FoP")
$newSheet = $wb.Worksheets.Add($null, $fop)
$newSheet.Name = "FiP"

# ---------------------------------------------------------------------------
# 2. Header row (row 1) — identical layout to the other preference sheets:
#    No. of responses | DIP | AIS | SEO | SA | UE | ACN | Sum || DIP | AIS |
#    SEO | SA | UE | ACN | Total
# ---------------------------------------------------------------------------
$newSheet.Range("A1").Value = "No. of responses"
$newSheet.Range("B1").Value = "DIP"
$newSheet.Range("C1").Value = "AIS"
$newSheet.Range("D1").Value = "SEO"
$newSheet.Range("E1").Value = "SA"
$newSheet.Range("F1").Value = "UE"
$newSheet.Range("G1").Value = "ACN"
$newSheet.Range("H1").Value = "Sum"

$newSheet.Range("J1").Value = "DIP"
$newSheet.Range("K1").Value = "AIS"
$newSheet.Range("L1").Value = "SEO"
$newSheet.Range("M1").Value = "SA"
$newSheet.Range("N1").Value = "UE"
$newSheet.Range("O1").Value = "ACN"
$newSheet.Range("P1").Value = "Total"

# ---------------------------------------------------------------------------
# 3. Raw input data (column A = number of responses, J:O = percentages)
# ---------------------------------------------------------------------------
$data = @(
    @{ Row = 2; A = 19; J = 10.5; K = 15.8; L = 15.8; M = 36.8; N = 5.3;  O = 15.8 },
    @{ Row = 3; A = 38; J = 7.9;  K = 26.3; L = 10.5; M = 34.2; N = 2.7;  O = 18.4 },
    @{ Row = 4; A = 54; J = 13;   K = 27.8; L = 6.9;  M = 31.5; N = 2.3;  O = 18.5 },
    @{ Row = 5; A = 63; J = 12.7; K = 25.4; L = 5.9;  M = 27;   N = 2;    O = 27   }
)

foreach ($d in $data) {
    $r = $d.Row

    $newSheet.Range("A$r").Value = $d.A
    $newSheet.Range("J$r").Value = $d.J
    $newSheet.Range("K$r").Value = $d.K
    $newSheet.Range("L$r").Value = $d.L
    $newSheet.Range("M$r").Value = $d.M
    $newSheet.Range("N$r").Value = $d.N
    $newSheet.Range("O$r").Value = $d.O

    # Computed columns (mirrors every other preference sheet)
    $newSheet.Range("B$r").Formula = "=A$r*J$r%"
    $newSheet.Range("C$r").Formula = "=A$r*K$r%"
    $newSheet.Range("D$r").Formula = "=A$r*L$r%"
    $newSheet.Range("E$r").Formula = "=A$r*M$r%"
    $newSheet.Range("F$r").Formula = "=A$r*N$r%"
    $newSheet.Range("G$r").Formula = "=A$r*O$r%"
    $newSheet.Range("H$r").Formula = "=SUM(B$r`:G$r)"
    $newSheet.Range("P$r").Formula = "=SUM(J$r`:O$r)"

    # Number format ("0") applied to the computed columns, matching the
    # other preference sheets (style index reused, no new style created).
    $newSheet.Range("B$r`:H$r").NumberFormat = "0"
    $newSheet.Range("I$r").NumberFormat = "0"
}

# ---------------------------------------------------------------------------
# 4. Sheet views / selection bookkeeping.
#    The newly added sheet becomes the active / selected tab, while FoP
#    loses its tabSelected flag and its selection resets to the full range.
# ---------------------------------------------------------------------------
$fop.Activate() | Out-Null
$fop.Range("A1:P5").Select() | Out-Null

$newSheet.Activate() | Out-Null
$newSheet.Range("N7").Select() | Out-Null
